$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1, J1 matching the existing header style (copy format from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill I2:J69 with new data values
$arr = New-Object "object[,]" 68,2
$arr[0,0] = 8
$arr[0,1] = 8
$arr[1,0] = 8
$arr[1,1] = 9
$arr[2,0] = 7
$arr[2,1] = 7
$arr[3,0] = 8
$arr[3,1] = 8
$arr[4,0] = 9
$arr[4,1] = 9
$arr[5,0] = 9
$arr[5,1] = 9
$arr[6,0] = 9
$arr[6,1] = 9
$arr[7,0] = 8
$arr[7,1] = 9
$arr[8,0] = 8
$arr[8,1] = 9
$arr[9,0] = 9
$arr[9,1] = 9
$arr[10,0] = 8
$arr[10,1] = 9
$arr[11,0] = 8
$arr[11,1] = 8
$arr[12,0] = 8
$arr[12,1] = 9
$arr[13,0] = 8
$arr[13,1] = 8
$arr[14,0] = 9
$arr[14,1] = 9
$arr[15,0] = 7
$arr[15,1] = 8
$arr[16,0] = 7
$arr[16,1] = 7
$arr[17,0] = 8
$arr[17,1] = 8
$arr[18,0] = 6
$arr[18,1] = 7
$arr[19,0] = 8
$arr[19,1] = 8
$arr[20,0] = 8
$arr[20,1] = 8
$arr[21,0] = 8
$arr[21,1] = 8
$arr[22,0] = 9
$arr[22,1] = 9
$arr[23,0] = 6
$arr[23,1] = 6
$arr[24,0] = 8
$arr[24,1] = 8
$arr[25,0] = 7
$arr[25,1] = 7
$arr[26,0] = 6
$arr[26,1] = 7
$arr[27,0] = 7
$arr[27,1] = 7
$arr[28,0] = 7
$arr[28,1] = 8
$arr[29,0] = 6
$arr[29,1] = 7
$arr[30,0] = 8
$arr[30,1] = 8
$arr[31,0] = 5
$arr[31,1] = 6
$arr[32,0] = 6
$arr[32,1] = 7
$arr[33,0] = 6
$arr[33,1] = 6
$arr[34,0] = 8
$arr[34,1] = 8
$arr[35,0] = 5
$arr[35,1] = 6
$arr[36,0] = 5
$arr[36,1] = 6
$arr[37,0] = 4
$arr[37,1] = 5
$arr[38,0] = 7
$arr[38,1] = 8
$arr[39,0] = 7
$arr[39,1] = 8
$arr[40,0] = 7
$arr[40,1] = 7
$arr[41,0] = 5
$arr[41,1] = 6
$arr[42,0] = 7
$arr[42,1] = 7
$arr[43,0] = 7
$arr[43,1] = 7
$arr[44,0] = 7
$arr[44,1] = 7
$arr[45,0] = 6
$arr[45,1] = 7
$arr[46,0] = 7
$arr[46,1] = 7
$arr[47,0] = 6
$arr[47,1] = 6
$arr[48,0] = 6
$arr[48,1] = 6
$arr[49,0] = 8
$arr[49,1] = 8
$arr[50,0] = 6
$arr[50,1] = 6
$arr[51,0] = 6
$arr[51,1] = 7
$arr[52,0] = 8
$arr[52,1] = 9
$arr[53,0] = 8
$arr[53,1] = 8
$arr[54,0] = 8
$arr[54,1] = 8
$arr[55,0] = 9
$arr[55,1] = 9
$arr[56,0] = 8
$arr[56,1] = 8
$arr[57,0] = 7
$arr[57,1] = 8
$arr[58,0] = 6
$arr[58,1] = 6
$arr[59,0] = 8
$arr[59,1] = 8
$arr[60,0] = 8
$arr[60,1] = 8
$arr[61,0] = 5
$arr[61,1] = 6
$arr[62,0] = 4
$arr[62,1] = 4
$arr[63,0] = 4
$arr[63,1] = 5
$arr[64,0] = 6
$arr[64,1] = 6
$arr[65,0] = 4
$arr[65,1] = 4
$arr[66,0] = 5
$arr[66,1] = 5
$arr[67,0] = 5
$arr[67,1] = 5
$ws.Range("I2:J69").Value = $arr
